# Color_calc.xlsx - "icons (only) in the bottom navigation"
#
# Sheet1 is an icon-size calculator for Android density buckets:
#   column C=ldpi, D=mdpi, E=hdpi, F=xhdpi, G=xxhdpi, H=xxxhdpi
#   row 10=0.75x, row 11=1x (base), row 12=1.5x, row 13=2x, row 14=3x, row 15=4x
# Each density column has exactly one directly-entered base size (row 11,
# except ldpi which is entered at row 10); every other cell in that column
# is a formula derived from it.
#
# This edit:
#   - lowers the mdpi (column D) base icon size from 32 to 24 (px), which
#     cascades through the 0.75x/1.5x/2x/3x/4x formulas in that column
#   - clears the xxhdpi (column G) entry at G14 (was 168), so the formulas
#     that key off it fall back to 0
#   - moves the active selection to D15 (the last recalculated mdpi cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Base mdpi icon size: 32 -> 24
$ws.Range("D11").Value = 24

# Clear the xxhdpi entry that used to drive column G
$ws.Range("G14").ClearContents()

# Leave the selection where the author left it after the edit
$ws.Range("D15").Select()
